# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# worksheets to reflect the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> new F-column value
$updates = @{
    2  = 8886
    3  = 8358
    4  = 152
    5  = 171
    6  = 219
    8  = 767
    9  = 224
    10 = 5577
    11 = 14
    13 = 91
    17 = 172
    18 = 237
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
